# Chỉnh giao diện code
# Insert 4 new receipt rows above the existing one and refresh the numbering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 4 new rows right before the current (only) data row.
$ws.Rows("4:7").Insert()
# Insert() copies formatting down from the row above (the bold header) -
# reset the newly inserted rows back to the default "Normal" style used
# by the rest of the data rows.
$ws.Range("A4:F7").Style = "Normal"

# Row 4
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "PT002"
$ws.Range("C4").Value = "Trần Nguyễn Yến Nhi"
$ws.Range("D4").Value = "Cẩm nang chăm sóc sức khỏe"
$ws.Range("E4").Value = 17
$ws.Range("F4").Value = 34000

# Row 5
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "PT007"
$ws.Range("C5").Value = "Trần Nguyễn Yến Nhi"
$ws.Range("D5").Value = "Toán học và ứng dụng"
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 18000

# Row 6
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "PT009"
$ws.Range("C6").Value = "Nguyễn Thanh Hưng"
$ws.Range("D6").Value = "Giáo trình Hệ điều hành"
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 8000

# Row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "PT010"
$ws.Range("C7").Value = "Lê Thị Ngọc Ánh"
$ws.Range("D7").Value = "Ngôn ngữ lập trình C#"
$ws.Range("E7").Value = 23
$ws.Range("F7").Value = 46000

# Row 8 already holds the original record (PT011 / Trần Lê Tuyết Mai /
# Đại số tuyến tính / 33 / 66000) - it was pushed down automatically by
# the row insert above, only its running number needs to stay 5.
$ws.Range("A8").Value = 5

# Widen the name / book-title columns so the longer Vietnamese text fits.
$ws.Columns("C:D").AutoFit()
